$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column C
$ws.Range("C1").Value = "Sex"

# Fill in Sex values for rows 2-15
$sexValues = @("M", "M", "F", "F", "NA", "F", "M", "M", "F", "M", "F", "F", "M", "M")
for ($i = 0; $i -lt $sexValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $sexValues[$i]
}

# Update selection to C14 as in the diff
$ws.Range("C14").Select()
